# Fixed Unique Jobs not printing.
# Added Assign and Write Cell activities to If CompanySize Else statement.
#
# The "UniqueJobs" column (G) was not being written for the rows that fall
# into the CompanySize "Else" branch (rows 3, 5 and 7). Backfill those cells
# with the same "N/A" placeholder used elsewhere in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CareerBuilder")

$ws.Range("G3").Value = "N/A"
$ws.Range("G5").Value = "N/A"
$ws.Range("G7").Value = "N/A"

# Reflect the updated view state: zoomed out a bit and the selection
# extended to cover the now-complete data range.
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("A2:L7").Select()
